$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the old used range (rows 1-9, cols A-D) before rewriting
$ws.Range("A1:D9").Clear()

# Row 1: header "shot" + 4 shot numbers
$ws.Range("A1").Value = "shot"
$ws.Range("B1").Value = 2661
$ws.Range("C1").Value = 2659
$ws.Range("D1").Value = 2615
$ws.Range("E1").Value = 2614

# Row 2: "sch frame 1" + values
$ws.Range("A2").Value = "sch frame 1"
$ws.Range("B2").Value = 67.6
$ws.Range("C2").Value = 51.9
$ws.Range("D2").Value = 47.5
$ws.Range("E2").Value = 85.7

# Row 3: "sch frame 2" + values
$ws.Range("A3").Value = "sch frame 2"
$ws.Range("B3").Value = 97.6
$ws.Range("C3").Value = 81.9
$ws.Range("D3").Value = 77.5

# Update selection to match target (E3 selected)
$ws.Range("E3").Select()
